# "filtered by banks and by concept key"
#
# Sheet1: add a header row (cve_concepto / descripcion) above the concept
# table, and close up the blank row that used to separate A3:B9 from
# A11:B23 so the table becomes one contiguous block right under the new
# header. The data-validation list (column B) follows the rows it was
# attached to automatically.
#
# "bancos" sheet: the banks column (A2:A12) is selected/highlighted - this
# is the list later used to build the "Bancos" filter.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# 1) Close the blank separator row (old row 10) so the two blocks become
#    contiguous (rows 3-9 then 11-23 -> rows 3-22).
$ws1.Rows.Item(10).Delete()

# 2) Rows 1-2 are currently unused padding above the table (table used to
#    start at row 3) - remove them so the table starts right at row 1.
$ws1.Range("A1:A2").EntireRow.Delete()

# 3) Insert the new header row above the table and fill it in.
$ws1.Rows.Item(1).Insert()
$ws1.Range("A1").Value = "cve_concepto"
$ws1.Range("B1").Value = "descripcion"

# "bancos" sheet: select the bank names column.
$ws3 = $wb.Worksheets.Item("bancos")
$ws3.Range("A2:A12").Select() | Out-Null
